# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (with fresh fund-holding data) between
# the existing "总计" (totals) sheet and the existing "2022-Q2" sheet, and
# appends a corresponding "2022-Q3" row to the "总计" summary sheet (pushing
# the old "2022-Q2" summary row down to row 3).

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)      # "总计"
$q2    = $wb.Worksheets.Item(2)      # "2022-Q2" (existing)

# ---------------------------------------------------------------------
# 1) "总计": shift the existing 2022-Q2 summary row down to row 3 (keeping
#    its formatting), then overwrite row 2 with the new 2022-Q3 numbers.
# ---------------------------------------------------------------------
$total.Range("A2:D2").Copy($total.Range("A3:D3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.44

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.42

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q3" sheet by copying "总计" (so it inherits the
#    same header/row-label styling) and placing it right after "总计",
#    i.e. before the existing "2022-Q2" sheet.
# ---------------------------------------------------------------------
$total.Copy($null, $total)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Extend the bold/bordered header style across E1:H1 and the row-label
# style down A3:A4, reusing the styles already present on B1 / A2.
$q3.Range("B1").Copy($q3.Range("E1:H1"))
$q3.Range("A2").Copy($q3.Range("A3:A4"))

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# The B column (fund codes) and D:G columns (numeric-looking figures) are
# stored as TEXT in the source data, so force a text number-format before
# assigning them (otherwise e.g. "000041" would become the number 41).
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

# Row 2
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "000041"
$q3.Range("C2").Value = "华夏全球精选股票（QDII）"
$q3.Range("D2").Value = "18.51"
$q3.Range("E2").Value = "89.41"
$q3.Range("F2").Value = "2.08"
$q3.Range("G2").Value = "0.3850"
$q3.Range("H2").Value = 10

# Row 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "519601"
$q3.Range("C3").Value = "海富通中国海外精选混合（QDII）"
$q3.Range("D3").Value = "0.51"
$q3.Range("E3").Value = "73.52"
$q3.Range("F3").Value = "5.00"
$q3.Range("G3").Value = "0.0255"
$q3.Range("H3").Value = 3

# Row 4
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "519602"
$q3.Range("C4").Value = "海富通大中华精选混合（QDII）"
$q3.Range("D4").Value = "0.10"
$q3.Range("E4").Value = "87.37"
$q3.Range("F4").Value = "5.89"
$q3.Range("G4").Value = "0.0059"
$q3.Range("H4").Value = 3

# ---------------------------------------------------------------------
# 3) Leave the focus/active tab on "总计", matching the original workbook.
# ---------------------------------------------------------------------
$total.Activate()
$total.Range("A1").Select()
